# Append a new row (row 72) to each of the 4 worksheets, duplicating the
# last existing row (row 71) except for the timestamp in column A, which
# advances to the next day's reading.

$wb = $excel.ActiveWorkbook

$newRow = 72
$newTime = 45858.43453703704
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Per-sheet payload for columns B..I, taken verbatim from the existing row 71
# of each sheet (the new row duplicates everything except the time stamp).
$rowData = @(
    @{ B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x48"; E = "0x14"; F = 380; G = 759863127514710945038336.0; H = 328; I = 14 },
    @{ B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x48"; E = "0xe";  F = 380; G = 568432987514711010443264.0; H = 328; I = 14 },
    @{ B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x7A"; E = "0x7";  F = 130; G = 568631262647113970876416.0; H = 122; I = 7  },
    @{ B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x7A"; E = "0x3";  F = 130; G = 985046333984776009023488.0; H = 122; I = 3  }
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowData[$i - 1]

    $ws.Cells.Item($newRow, 1).Value = $newTime
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
